# Insert a new weekly price record as row 76, pushing all subsequent
# rows (old 76..176) down by one (to 77..177).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(76).Insert()

$ws.Cells.Item(76, 1).Value = 5
$ws.Cells.Item(76, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(76, 3).Value = "Maule"
$ws.Cells.Item(76, 4).Value = 44467
$ws.Cells.Item(76, 5).Value = 7
$ws.Cells.Item(76, 6).Value = 100114014
$ws.Cells.Item(76, 7).Value = "Betarraga"
$ws.Cells.Item(76, 8).Value = "Sin especificar"
$ws.Cells.Item(76, 9).Value = "Primera"
$ws.Cells.Item(76, 10).Value = 4000
$ws.Cells.Item(76, 11).Value = 650
$ws.Cells.Item(76, 12).Value = 650
$ws.Cells.Item(76, 13).Value = 650
$ws.Cells.Item(76, 14).Value = "$/paquete 5 unidades"
$ws.Cells.Item(76, 15).Value = "Región del Maule"
$ws.Cells.Item(76, 16).Value = 130
$ws.Cells.Item(76, 17).Value = 5
$ws.Cells.Item(76, 18).Value = "Hortaliza"
